$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range('D2').Style
$ws.Range('D2').Value = "'61.979.85"
$ws.Range('D2').Style = $style
$ws.Range('E2').Value = '  -0.65%  '

$style = $ws.Range('D3').Style
$ws.Range('D3').Value = "'3.427.37"
$ws.Range('D3').Style = $style
$ws.Range('E3').Value = '  -0.13%  '

$style = $ws.Range('D4').Style
$ws.Range('D4').Value = "'1.00"
$ws.Range('D4').Style = $style
$ws.Range('E4').Value = '  +0.00%  '

$style = $ws.Range('D5').Style
$ws.Range('D5').Value = "'409.27"
$ws.Range('D5').Style = $style
$ws.Range('E5').Value = '  +0.47%  '

$style = $ws.Range('D6').Style
$ws.Range('D6').Value = "'128.62"
$ws.Range('D6').Style = $style
$ws.Range('E6').Value = '  -3.47%  '

$ws.Range('E7').Value = '  +6.12%  '

$ws.Range('E8').Value = '  -0.06%  '

$ws.Range('E9').Value = '  +6.68%  '

$ws.Range('E10').Value = '  +2.63%  '

$style = $ws.Range('D11').Style
$ws.Range('D11').Value = "'42.76"
$ws.Range('D11').Style = $style
$ws.Range('E11').Value = '  +1.90%  '

$ws.Range('E12').Value = '  +47.05%  '

$ws.Range('E13').Value = '  +8.62%  '

$ws.Range('E14').Value = '  -0.17%  '

$style = $ws.Range('D15').Style
$ws.Range('D15').Value = "'21.43"
$ws.Range('D15').Style = $style
$ws.Range('E15').Value = '  +7.69%  '

$style = $ws.Range('D16').Style
$ws.Range('D16').Value = "'3.963.41"
$ws.Range('D16').Style = $style
$ws.Range('E16').Value = '  -0.10%  '

$style = $ws.Range('D17').Style
$ws.Range('D17').Value = "'3.367.70"
$ws.Range('D17').Style = $style
$ws.Range('E17').Value = '  -2.84%  '

$style = $ws.Range('D18').Style
$ws.Range('D18').Value = "'12.50"
$ws.Range('D18').Style = $style
$ws.Range('E18').Value = '  +7.72%  '

$ws.Range('E19').Value = '  +6.82%  '

$style = $ws.Range('D20').Style
$ws.Range('D20').Value = "'61.923.03"
$ws.Range('D20').Style = $style
$ws.Range('E20').Value = '  -0.48%  '

$style = $ws.Range('D21').Style
$ws.Range('D21').Value = "'453.45"
$ws.Range('D21').Style = $style
$ws.Range('E21').Value = '  +45.10%  '

$style = $ws.Range('D22').Style
$ws.Range('D22').Value = "'91.97"
$ws.Range('D22').Style = $style
$ws.Range('E22').Value = '  +9.00%  '

$ws.Range('E23').Value = '  +0.96%  '

$style = $ws.Range('D24').Style
$ws.Range('D24').Value = "'13.01"
$ws.Range('D24').Style = $style
$ws.Range('E24').Value = '  +1.60%  '

$style = $ws.Range('D25').Style
$ws.Range('D25').Value = "'3.23"
$ws.Range('D25').Style = $style
$ws.Range('E25').Value = '  +1.95%  '

$style = $ws.Range('D26').Style
$ws.Range('D26').Value = "'32.93"
$ws.Range('D26').Style = $style
$ws.Range('E26').Value = '  +10.65%  '

$style = $ws.Range('D27').Style
$ws.Range('D27').Value = "'8.80"
$ws.Range('D27').Style = $style
$ws.Range('E27').Value = '  +7.81%  '

$ws.Range('E28').Value = '  +0.44%  '

$ws.Range('E29').Value = '  -2.18%  '

$ws.Range('E30').Value = '  -2.64%  '

$style = $ws.Range('D31').Style
$ws.Range('D31').Value = "'12.00"
$ws.Range('D31').Style = $style
$ws.Range('E31').Value = '  +5.70%  '

$ws.Range('E32').Value = '  -0.77%  '

$style = $ws.Range('D33').Style
$ws.Range('D33').Value = "'43.10"
$ws.Range('D33').Style = $style
$ws.Range('E33').Value = '  -1.74%  '

$ws.Range('E34').Value = '  -0.61%  '

$ws.Range('E36').Value = '  +3.05%  '

$style = $ws.Range('D37').Style
$ws.Range('D37').Value = "'54.39"
$ws.Range('D37').Style = $style
$ws.Range('E37').Value = '  +5.29%  '

$ws.Range('E38').Value = '  +0.12%  '

$ws.Range('E39').Value = '  +1.26%  '

$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$style = $ws.Range('D40').Style
$ws.Range('D40').Value = "'0.135"
$ws.Range('D40').Style = $style
$ws.Range('E40').Value = '  +7.46%  '

$ws.Range('B41').Value = 'TheGraph'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$style = $ws.Range('D41').Style
$ws.Range('D41').Value = "'0.323"
$ws.Range('D41').Style = $style
$ws.Range('E41').Value = '  +2.18%  '

$ws.Range('E42').Value = '  -2.89%  '

$style = $ws.Range('D43').Style
$ws.Range('D43').Value = "'142.20"
$ws.Range('D43').Style = $style
$ws.Range('E43').Value = '  +0.30%  '

$ws.Range('E44').Value = '  +8.33%  '

$ws.Range('E45').Value = '  +0.98%  '

$style = $ws.Range('D46').Style
$ws.Range('D46').Value = "'2.51"
$ws.Range('D46').Style = $style
$ws.Range('E46').Value = '  +13.13%  '

$style = $ws.Range('D47').Style
$ws.Range('D47').Value = "'16.68"
$ws.Range('D47').Style = $style
$ws.Range('E47').Value = '  -0.54%  '

$style = $ws.Range('D48').Style
$ws.Range('D48').Value = "'22.37"
$ws.Range('D48').Style = $style
$ws.Range('E48').Value = '  +4.66%  '

$ws.Range('E49').Value = '  +9.98%  '

$style = $ws.Range('D50').Style
$ws.Range('D50').Value = "'3.773.94"
$ws.Range('D50').Style = $style
$ws.Range('E50').Value = '  +0.14%  '

$style = $ws.Range('D51').Style
$ws.Range('D51').Value = "'0.138"
$ws.Range('D51').Style = $style
$ws.Range('E51').Value = '  +15.60%  '
